$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '36.490.62'
$ws.Range('E2').Value2 = '  +0.09%  '
$ws.Range('D3').Value2 = '1.938.41'
$ws.Range('E3').Value2 = '  -2.03%  '
$ws.Range('E4').Value2 = '  -0.05%  '
$ws.Range('D5').Value2 = '''242.42'
$ws.Range('E5').Value2 = '  -1.30%  '
$ws.Range('D6').Value2 = '''0.609'
$ws.Range('E6').Value2 = '  -2.71%  '
$ws.Range('E7').Value2 = '  -0.06%  '
$ws.Range('D8').Value2 = '''56.88'
$ws.Range('E8').Value2 = '  -3.75%  '
$ws.Range('D9').Value2 = '''0.359'
$ws.Range('E9').Value2 = '  -3.94%  '
$ws.Range('D10').Value2 = '''0.0847'
$ws.Range('E10').Value2 = '  +0.60%  '
$ws.Range('E11').Value2 = '  -1.60%  '
$ws.Range('D12').Value2 = '2.216.85'
$ws.Range('E12').Value2 = '  -2.24%  '
$ws.Range('D13').Value2 = '''21.28'
$ws.Range('E13').Value2 = '  -6.60%  '
$ws.Range('D14').Value2 = '''0.807'
$ws.Range('E14').Value2 = '  -5.96%  '
$ws.Range('D15').Value2 = '''13.41'
$ws.Range('E15').Value2 = '  -3.58%  '
$ws.Range('D16').Value2 = '''5.13'
$ws.Range('E16').Value2 = '  -5.99%  '
$ws.Range('D17').Value2 = '1.940.72'
$ws.Range('E17').Value2 = '  -2.14%  '
$ws.Range('D18').Value2 = '36.373.75'
$ws.Range('E18').Value2 = '  +0.07%  '
$ws.Range('D19').Value2 = '''69.09'
$ws.Range('E19').Value2 = '  -1.98%  '
$ws.Range('D20').Value2 = '0.0₃0861'
$ws.Range('E20').Value2 = '  -2.26%  '
$ws.Range('D21').Value2 = '''226.73'
$ws.Range('E21').Value2 = '  -3.13%  '
$ws.Range('D22').Value2 = '''4.97'
$ws.Range('E22').Value2 = '  -5.69%  '
$ws.Range('E23').Value2 = '  -0.20%  '
$ws.Range('D24').Value2 = '''2.32'
$ws.Range('E24').Value2 = '  -7.89%  '
$ws.Range('E25').Value2 = '  -1.07%  '
$ws.Range('D26').Value2 = '''9.11'
$ws.Range('E26').Value2 = '  -7.76%  '
$ws.Range('D27').Value2 = '''160.92'
$ws.Range('E27').Value2 = '  -1.67%  '
$ws.Range('D28').Value2 = '''0.135'
$ws.Range('E28').Value2 = '  +1.34%  '
$ws.Range('D29').Value2 = '''19.23'
$ws.Range('E29').Value2 = '  -3.04%  '
$ws.Range('D30').Value2 = '''0.118'
$ws.Range('E30').Value2 = '  -1.69%  '
$ws.Range('E31').Value2 = '  -5.89%  '
$ws.Range('D32').Value2 = '''4.54'
$ws.Range('E32').Value2 = '  -6.84%  '
$ws.Range('D33').Value2 = '''0.0620'
$ws.Range('E33').Value2 = '  -7.04%  '
$ws.Range('D34').Value2 = '''4.15'
$ws.Range('E34').Value2 = '  -6.89%  '
$ws.Range('E35').Value2 = '  +0.02%  '
$ws.Range('B36').Value2 = 'THORChain'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').Value2 = '''5.99'
$ws.Range('E36').Value2 = '  -2.36%  '
$ws.Range('B37').Value2 = 'WEMIXToken'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value2 = '''1.80'
$ws.Range('E37').Value2 = '  -0.97%  '
$ws.Range('D38').Value2 = '''2.16'
$ws.Range('E38').Value2 = '  -2.99%  '
$ws.Range('D39').Value2 = '''3.06'
$ws.Range('E39').Value2 = '  +3.55%  '
$ws.Range('D40').Value2 = '''0.0995'
$ws.Range('E40').Value2 = '  +2.95%  '
$ws.Range('E41').Value2 = '  -0.29%  '
$ws.Range('D42').Value2 = '''0.0209'
$ws.Range('E42').Value2 = '  -2.42%  '
$ws.Range('E43').Value2 = '  -5.91%  '
$ws.Range('D44').Value2 = '''15.59'
$ws.Range('E44').Value2 = '  -3.24%  '
$ws.Range('D45').Value2 = '1.338.21'
$ws.Range('E45').Value2 = '  -2.52%  '
$ws.Range('D46').Value2 = '''1.02'
$ws.Range('E46').Value2 = '  -6.46%  '
$ws.Range('D47').Value2 = '''85.63'
$ws.Range('E47').Value2 = '  -6.41%  '
$ws.Range('D48').Value2 = '''7.08'
$ws.Range('E48').Value2 = '  -4.66%  '
$ws.Range('D49').Value2 = '''2.83'
$ws.Range('E49').Value2 = '  -0.38%  '
$ws.Range('D50').Value2 = '2.108.54'
$ws.Range('E50').Value2 = '  -2.24%  '
$ws.Range('D51').Value2 = '''43.38'
$ws.Range('E51').Value2 = '  -4.29%  '
